# Round the numeric data (columns B:E, rows 2:13) to the nearest integer.
# The source values were stored with full floating point precision; the
# target output stores them as plain rounded integers (values only, the
# cell number format stays General).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 13
$firstCol = 2   # B
$lastCol  = 5   # E

for ($r = 2; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = [double]$cell.Value2
        $cell.Value2 = [math]::Floor($v + 0.5)
    }
}
